# Populate the (previously empty) Sheet1 with the single status/error
# message produced by the config-file detection pass, matching the
# "added support for configuration files detecting ... unconditional
# looping over data and generating results" export behaviour: one row,
# one shared-string cell, with the column auto-sized to the text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$message = "error: unsupported file:c:\gitrep\xmlparse\dell\2018-12-6_85625_export.xml"

$ws.Range("A1").Value = $message

# Size column A to fit the message text (character-width units, as
# Excel stores them in the sheet's <cols> definition).
$ws.Columns.Item(1).ColumnWidth = 74.7109375
